$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'156"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'343600.00"
$ws.Range("D2").Style = "Normal"

$ws.Range("C3").Value = "'854"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'2327070.65"
$ws.Range("D3").Style = "Normal"

$ws.Range("C4").Value = "'350"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1236968.92"
$ws.Range("D4").Style = "Normal"

$ws.Range("C5").Value = "'91"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'367982.09"
$ws.Range("D5").Style = "Normal"

$ws.Range("C15").Value = "'91"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'224152.38"
$ws.Range("D15").Style = "Normal"

$ws.Range("C16").Value = "'400"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1185147.19"
$ws.Range("D16").Style = "Normal"

$ws.Range("C18").Value = "'42"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'196045.00"
$ws.Range("D18").Style = "Normal"

$ws.Range("C34").Value = "'468"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1325978.53"
$ws.Range("D34").Style = "Normal"

$ws.Range("C35").Value = "'186"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'820539.11"
$ws.Range("D35").Style = "Normal"

$ws.Range("C38").Value = "'14"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'28000.00"
$ws.Range("D38").Style = "Normal"

$ws.Range("C39").Value = "'29"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'72330.00"
$ws.Range("D39").Style = "Normal"

$ws.Range("C40").Value = "'144"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'363579.00"
$ws.Range("D40").Style = "Normal"

$ws.Range("C41").Value = "'74"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'250900.00"
$ws.Range("D41").Style = "Normal"

$ws.Range("C42").Value = "'18"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'65995.14"
$ws.Range("D42").Style = "Normal"

$ws.Range("C43").Value = "'10"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'44500.00"
$ws.Range("D43").Style = "Normal"

$ws.Range("C44").Value = "'32"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'68905.00"
$ws.Range("D44").Style = "Normal"

$ws.Range("C45").Value = "'20"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'88621.84"
$ws.Range("D45").Style = "Normal"

$ws.Range("C47").Value = "'30"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'173937.00"
$ws.Range("D47").Style = "Normal"

$ws.Range("C48").Value = "'22"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'143697.00"
$ws.Range("D48").Style = "Normal"
